$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# --- New header labels (row 1) ---
# Order chosen so the shared-strings table gets populated in the same
# sequence as the target file (ia, Vr, Ra, Rr, "Tm (v^2 /ra)/w", "Tm=va*ia/wm").
$ws.Range("G1").Value = "ia"
$ws.Range("F1").Value = "Vr"
$ws.Range("N2").Value = "Ra"
$ws.Range("N1").Value = "Rr"
$ws.Range("J1").Value = "Tm (v^2 /ra)/w"
$ws.Range("I1").Value = "Tm=va*ia/wm"

# --- Constants used by the formulas below ---
$ws.Range("O1").Value = 3.8
$ws.Range("O2").Value = 7.9

# --- New measured value ---
$ws.Range("F5").Value = 1.1200000000000001

# --- Column G: Vr/O1 ratio, one independent formula per row ---
for ($r = 2; $r -le 8; $r++) {
    $ws.Range("G$r").Formula = "=F$r/`$O`$1"
}

# --- Column I: torque from measured Va*ia/w (row 5 only) ---
$ws.Range("I5").Formula = "=B5*G5/C5"

# --- Column J: torque from (Va^2)/(Ra*w) ---
#   J2 is its own formula; J3:J8 are filled together so Excel stores them
#   as one shared-formula group, matching the source file.
$ws.Range("J2").Formula = "=(B2*B2)/(`$O`$2*C2)"
$ws.Range("J3:J8").Formula = "=(B3*B3)/(`$O`$2*C3)"

# --- Number formats ---
$ws.Range("I5").NumberFormat = "0.00000"
$ws.Range("G5").NumberFormat = "0.000"

# --- Column widths (auto-fit like Excel would do after typing the headers) ---
$ws.Columns.Item(9).AutoFit()
$ws.Columns.Item(10).AutoFit()

# --- Selection matches the author's final cursor position ---
$ws.Range("I5").Select()
